# msz - first smoke test is running
#
# Fills in the new "Goto product page" / "ProductData" / "SendQuote" rows of
# the Tabelle1 test-script sheet and repositions the screenshot picture that
# sits below the table so it keeps its two-row gap under the (now longer)
# table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# The table originally only had rows 6-8 (FillPageVehicleData / FillPageInsurantData
# / "Button Next from Page VehicleData"); the edit inserts extra rows among them, so
# every existing row 6-8 is first cleared completely and then rewritten at its new
# position to avoid leaving stale cells behind in columns that are no longer used.
$ws.Range("A6:H8").ClearContents()

# --- Row 6 gets new content ("Goto product page" / <SELECT>) -------------
$ws.Range("A6").Value = "Goto product page"
$ws.Range("B6").Value = "<SET>"
$ws.Range("E6").Value = "<SELECT>"

# --- Row 7: FillPageVehicleData -> FillPage (moved down from the old row 6)
$ws.Range("A7").Value = "102_VehicleInsuranceAutomobile_001_SmokeTest_FillPageVehicleData"
$ws.Range("B7").Value = "<SET>"
$ws.Range("C7").Value = "102_VehicleInsuranceAutomobile_001_SmokeTest_FillPage"
$ws.Range("H7").Value = "<NOP>"

# --- Row 8: FillPageInsurantData (moved down from the old row 7) ---------
$ws.Range("A8").Value = "102_VehicleInsuranceAutomobile_001_SmokeTest_FillPageInsurantData"
$ws.Range("B8").Value = "<SET>"
$ws.Range("D8").Value = "102_VehicleInsuranceAutomobile_001_SmokeTest_FillPage"
$ws.Range("H8").Value = "<NOP>"

# --- Row 9: new FillPageProductData block ---------------------------------
$ws.Range("A9").Value = "102_VehicleInsuranceAutomobile_001_SmokeTest_FillPageProductData"
$ws.Range("B9").Value = "<SET>"
$ws.Range("E9").Value = "102_VehicleInsuranceAutomobile_001_SmokeTest_FillPage"
$ws.Range("H9").Value = "<NOP>"

# --- Row 10: new FillPageSendQuote block ----------------------------------
$ws.Range("A10").Value = "102_VehicleInsuranceAutomobile_001_SmokeTest_FillPageSendQuote"
$ws.Range("B10").Value = "<SET>"
$ws.Range("G10").Value = "102_VehicleInsuranceAutomobile_001_SmokeTest_FillPage"
$ws.Range("H10").Value = "<NOP>"

# --- Row 11: Button Next from Page VehicleData (moved down from old row 8)
$ws.Range("A11").Value = "Button Next from Page VehicleData"
$ws.Range("B11").Value = "<SET>"
$ws.Range("C11").Value = "Button Next"
$ws.Range("H11").Value = "<NOP>"

# --- Rows 12-15: new "Choose <tier>" product rows -------------------------
$ws.Range("A12").Value = "Choose Silver"
$ws.Range("B12").Value = "<SET>"
$ws.Range("F12").Value = "Choose Silver"
$ws.Range("H12").Value = "<NOP>"

$ws.Range("A13").Value = "Choose Gold"
$ws.Range("B13").Value = "<SET>"
$ws.Range("F13").Value = "Choose Gold"
$ws.Range("H13").Value = "<NOP>"

$ws.Range("A14").Value = "Choose Platinum"
$ws.Range("B14").Value = "<SET>"
$ws.Range("F14").Value = "Choose Platinum"
$ws.Range("H14").Value = "<NOP>"

$ws.Range("A15").Value = "Choose Ultimate"
$ws.Range("B15").Value = "<SET>"
$ws.Range("F15").Value = "Choose Ultimate"
$ws.Range("H15").Value = "<NOP>"

# --- Row 16: new "Send Quote" row -----------------------------------------
$ws.Range("A16").Value = "Send Quote - Button Main Page"
$ws.Range("B16").Value = "<SET>"
$ws.Range("G16").Value = "Button Main Page"
$ws.Range("H16").Value = "<NOP>"

# --- Column widths: C:E and G now show the same "50, best-fit" width that
#     C:D already had (G used to be narrower before it held long text) -----
$ws.Columns("C").ColumnWidth = 49.17
$ws.Columns("D").ColumnWidth = 49.17
$ws.Columns("E").ColumnWidth = 49.17
$ws.Columns("G").ColumnWidth = 49.17

# --- Selection moved to C1:G2 (header block) ------------------------------
$ws.Range("C1:G2").Select()

# --- Move the screenshot picture down so it keeps its gap below the table -
$shp = $ws.Shapes.Item(1)
$shp.Top = $shp.Top + 111.6
